# Repurpose the "Werkt mee" (Ja/Nee) column on the "Data Sander" sheet into a
# new "Kwaliteit data" column that grades each company's data quality as
# Goud / Zilver / Brons.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Sander")

# Header
$ws.Range("F1").Value = "Kwaliteit data"

# Body - map old "Ja" -> "Zilver", old "Nee" -> "Brons", with Nippon (row 11)
# upgraded to the top "Goud" tier.
$ws.Range("F2").Value  = "Zilver"
$ws.Range("F3").Value  = "Zilver"
$ws.Range("F4").Value  = "Brons"
$ws.Range("F5").Value  = "Brons"
$ws.Range("F6").Value  = "Brons"
$ws.Range("F7").Value  = "Zilver"
$ws.Range("F8").Value  = "Zilver"
$ws.Range("F9").Value  = "Brons"
$ws.Range("F10").Value = "Zilver"
$ws.Range("F11").Value = "Goud"
$ws.Range("F12").Value = "Brons"
$ws.Range("F13").Value = "Zilver"
$ws.Range("F14").Value = "Zilver"
$ws.Range("F15").Value = "Brons"
$ws.Range("F16").Value = "Brons"
$ws.Range("F17").Value = "Brons"
$ws.Range("F18").Value = "Brons"
$ws.Range("F19").Value = "Brons"
$ws.Range("F20").Value = "Brons"
$ws.Range("F21").Value = "Brons"
$ws.Range("F22").Value = "Brons"
$ws.Range("F23").Value = "Brons"
$ws.Range("F24").Value = "Brons"
$ws.Range("F25").Value = "Brons"
$ws.Range("F26").Value = "Brons"
$ws.Range("F27").Value = "Brons"
$ws.Range("F28").Value = "Brons"
$ws.Range("F29").Value = "Brons"
$ws.Range("F30").Value = "Brons"
$ws.Range("F31").Value = "Brons"

# Leave the cursor where the author left it when they saved.
$ws.Range("F12").Select() | Out-Null
